# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

# Update OFF sheet (row labeled "R") with new cumulative stats
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 227
$wsOff.Range("C3").Value = 159
$wsOff.Range("D3").Value = 56
$wsOff.Range("E3").Value = 26
$wsOff.Range("F3").Value = 2
$wsOff.Range("G3").Value = 7

# Update DEF sheet (row labeled "R") with new cumulative stats
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 319
$wsDef.Range("C3").Value = 236
$wsDef.Range("D3").Value = 67
$wsDef.Range("E3").Value = 37
